$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2461.5217
$ws.Range("J17").Value = 2513.7778
$ws.Range("L17").Value = 7541.3334
$ws.Range("N17").Value = -7877.3334
# Row 33
$ws.Range("H33").Value = 5551934.5
$ws.Range("I33").Value = 9993304
$ws.Range("K33").Value = 9993304
$ws.Range("M33").Value = -9993075
# Row 40
$ws.Range("H40").Value = 2954.077
$ws.Range("I40").Value = 7000
$ws.Range("J40").Value = 1740.3
$ws.Range("K40").Value = 7000
$ws.Range("L40").Value = 1740.3
$ws.Range("M40").Value = -6825
$ws.Range("N40").Value = -2090.3
# Row 55
$ws.Range("H55").Value = 471.69232
$ws.Range("I55").Value = 250
$ws.Range("J55").Value = 570.2222
$ws.Range("K55").Value = 250
$ws.Range("L55").Value = 570.2222
$ws.Range("M55").Value = -36
$ws.Range("N55").Value = -998.2222
# Row 76
$ws.Range("H76").Value = 3716.6667
$ws.Range("I76").Value = 3462.5
$ws.Range("J76").Value = 5750
$ws.Range("K76").Value = 3462.5
$ws.Range("L76").Value = 5750
$ws.Range("M76").Value = -3147.5
$ws.Range("N76").Value = -6380
# Row 79
$ws.Range("H79").Value = 3716.6667
$ws.Range("I79").Value = 3462.5
$ws.Range("J79").Value = 5750
$ws.Range("K79").Value = 3462.5
$ws.Range("L79").Value = 5750
$ws.Range("M79").Value = -2370.5
$ws.Range("N79").Value = -7934
# Row 112
$ws.Range("H112").Value = 5584.647
$ws.Range("J112").Value = 1629.2667
$ws.Range("L112").Value = 4887.800099999999
$ws.Range("N112").Value = -7103.800099999999
# Row 129
$ws.Range("H129").Value = 942.2166999999999
$ws.Range("J129").Value = 998.0909
$ws.Range("L129").Value = 2994.2727
$ws.Range("N129").Value = -12994.2727
# Row 137
$ws.Range("H137").Value = 727668.0600000001
$ws.Range("I137").Value = 3472.2222
$ws.Range("J137").Value = 1193222.6
$ws.Range("K137").Value = 10416.6666
$ws.Range("L137").Value = 3579667.8
$ws.Range("M137").Value = -7866.6666
$ws.Range("N137").Value = -3584767.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 2008
$ws.Range("I19").Value = 2008
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2008
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1779
$ws.Range("N19").ClearContents()
# Row 32
$ws.Range("H32").Value = 18859.629
$ws.Range("I32").Value = 20254.803
$ws.Range("J32").Value = 5838
$ws.Range("K32").Value = 20254.803
$ws.Range("L32").Value = 5838
$ws.Range("M32").Value = -19967.803
$ws.Range("N32").Value = -6412
# Row 63
$ws.Range("H63").Value = 3493.5715
$ws.Range("I63").Value = 2991
$ws.Range("J63").Value = 4750
$ws.Range("K63").Value = 2991
$ws.Range("L63").Value = 4750
$ws.Range("M63").Value = -2305
$ws.Range("N63").Value = -6122
# Row 66
$ws.Range("H66").Value = 3493.5715
$ws.Range("I66").Value = 2991
$ws.Range("J66").Value = 4750
$ws.Range("K66").Value = 14955
$ws.Range("L66").Value = 23750
$ws.Range("M66").Value = -11523
$ws.Range("N66").Value = -30614
# Row 74
$ws.Range("H74").Value = 4670.2812
$ws.Range("I74").Value = 1571.6666
$ws.Range("J74").Value = 21402.8
$ws.Range("K74").Value = 1571.6666
$ws.Range("L74").Value = 21402.8
$ws.Range("M74").Value = -697.6666
$ws.Range("N74").Value = -23150.8
# Row 77
$ws.Range("H77").Value = 4670.2812
$ws.Range("I77").Value = 1571.6666
$ws.Range("J77").Value = 21402.8
$ws.Range("K77").Value = 7858.333000000001
$ws.Range("L77").Value = 107014
$ws.Range("M77").Value = -3490.333000000001
$ws.Range("N77").Value = -115750
# Row 122
$ws.Range("H122").Value = 5001714.5
$ws.Range("I122").Value = 1811.409
$ws.Range("J122").Value = 41667670
$ws.Range("K122").Value = 5434.227000000001
$ws.Range("L122").Value = 125003010
$ws.Range("M122").Value = -2984.227000000001
$ws.Range("N122").Value = -125007910
# Row 132
$ws.Range("H132").Value = 2450.434
$ws.Range("I132").Value = 2279.2368
$ws.Range("J132").Value = 2884.1333
$ws.Range("K132").Value = 6837.7104
$ws.Range("L132").Value = 8652.3999
$ws.Range("M132").Value = -4307.7104
$ws.Range("N132").Value = -13712.3999
# Row 139
$ws.Range("H139").Value = 47200
$ws.Range("I139").Value = 42800
$ws.Range("J139").Value = 49400
$ws.Range("K139").Value = 42800
$ws.Range("L139").Value = 49400
$ws.Range("M139").Value = -37660
$ws.Range("N139").Value = -59680

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 362.4
$ws.Range("I22").Value = 370
$ws.Range("J22").Value = 351
$ws.Range("K22").Value = 370
$ws.Range("L22").Value = 351
$ws.Range("M22").Value = -197
$ws.Range("N22").Value = -697
# Row 94
$ws.Range("H94").Value = 2066.5557
$ws.Range("I94").Value = 2066.5557
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2066.5557
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1615.5557
$ws.Range("N94").ClearContents()
# Row 107
$ws.Range("H107").Value = 1629.1111
$ws.Range("I107").Value = 868.5833
$ws.Range("J107").Value = 3150.1667
$ws.Range("K107").Value = 868.5833
$ws.Range("L107").Value = 3150.1667
$ws.Range("M107").Value = 1051.4167
$ws.Range("N107").Value = -6990.1667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
# Row 75
$ws.Range("H75").Value = 40260
$ws.Range("J75").Value = 40260
$ws.Range("L75").Value = 40260
$ws.Range("N75").Value = -42256
# Row 78
$ws.Range("H78").Value = 40260
$ws.Range("J78").Value = 40260
$ws.Range("L78").Value = 120780
$ws.Range("N78").Value = -130764
# Row 132
$ws.Range("H132").Value = 2389.45
$ws.Range("I132").Value = 2038.2667
$ws.Range("J132").Value = 3443
$ws.Range("K132").Value = 6114.800099999999
$ws.Range("L132").Value = 10329
$ws.Range("M132").Value = -3584.800099999999
$ws.Range("N132").Value = -15389
# Row 134
$ws.Range("H134").Value = 1902.8793
$ws.Range("I134").Value = 1504.4318
$ws.Range("J134").Value = 3155.1428
$ws.Range("K134").Value = 4513.2954
$ws.Range("L134").Value = 9465.428400000001
$ws.Range("M134").Value = -1978.2954
$ws.Range("N134").Value = -14535.4284

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 101
$ws.Range("H101").Value = 7029
$ws.Range("J101").Value = 7029
$ws.Range("L101").Value = 21087
$ws.Range("N101").Value = -25955

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
# Row 111
$ws.Range("H111").Value = 69000
$ws.Range("J111").Value = 69000
$ws.Range("L111").Value = 69000
$ws.Range("N111").Value = -75134
# Row 132
$ws.Range("H132").Value = 45515.77
$ws.Range("I132").Value = 93613.45
$ws.Range("J132").Value = 10244.134
$ws.Range("K132").Value = 280840.35
$ws.Range("L132").Value = 30732.402
$ws.Range("M132").Value = -278310.35
$ws.Range("N132").Value = -35792.402

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 4724.75
$ws.Range("I17").Value = 5450
$ws.Range("J17").Value = 3999.5
$ws.Range("K17").Value = 5450
$ws.Range("L17").Value = 3999.5
$ws.Range("M17").Value = -5280
$ws.Range("N17").Value = -4339.5
# Row 68
$ws.Range("H68").Value = 4068.3044
$ws.Range("I68").Value = 3150
$ws.Range("J68").Value = 4558.067
$ws.Range("K68").Value = 3150
$ws.Range("L68").Value = 4558.067
$ws.Range("M68").Value = -2401
$ws.Range("N68").Value = -6056.067
# Row 71
$ws.Range("H71").Value = 4068.3044
$ws.Range("I71").Value = 3150
$ws.Range("J71").Value = 4558.067
$ws.Range("K71").Value = 15750
$ws.Range("L71").Value = 22790.335
$ws.Range("M71").Value = -12006
$ws.Range("N71").Value = -30278.335
# Row 82
$ws.Range("H82").Value = 2318.4
$ws.Range("I82").Value = 1251
$ws.Range("J82").Value = 3030
$ws.Range("K82").Value = 1251
$ws.Range("L82").Value = 3030
$ws.Range("M82").Value = -890
$ws.Range("N82").Value = -3752
# Row 85
$ws.Range("H85").Value = 2318.4
$ws.Range("I85").Value = 1251
$ws.Range("J85").Value = 3030
$ws.Range("K85").Value = 1251
$ws.Range("L85").Value = 3030
$ws.Range("M85").Value = -3
$ws.Range("N85").Value = -5526
# Row 138
$ws.Range("H138").Value = 50599.5
$ws.Range("J138").Value = 50599.5
$ws.Range("L138").Value = 50599.5
$ws.Range("N138").Value = -60879.5
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Range("H86").Value = 22966.666
$ws.Range("J86").Value = 22966.666
$ws.Range("L86").Value = 22966.666
$ws.Range("N86").Value = -25212.666
# Row 89
$ws.Range("H89").Value = 22966.666
$ws.Range("J89").Value = 22966.666
$ws.Range("L89").Value = 114833.33
$ws.Range("N89").Value = -126065.33
# Row 107
$ws.Range("H107").Value = 2885
$ws.Range("I107").Value = 426.83334
$ws.Range("J107").Value = 4359.9
$ws.Range("K107").Value = 1280.50002
$ws.Range("L107").Value = 13079.7
$ws.Range("M107").Value = 639.4999800000001
$ws.Range("N107").Value = -16919.7
# Row 141
$ws.Range("H141").Value = 64400
$ws.Range("J141").Value = 64400
$ws.Range("L141").Value = 64400
$ws.Range("N141").Value = -74760
